$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily rows appended to the series (08-09-2021 and 09-09-2021),
# same TPM / facilidad permanente rates as the most recent prior rows.

# Write dates as text via a formula (text literal) then paste-special
# as values so the cells end up as plain shared-string text, matching
# how the existing "dd-mm-yyyy" entries are stored (avoids Excel's
# automatic date-recognition when assigning a literal string Value).
$ws.Range("A174").Formula = "=""08-09-2021"""
$ws.Range("A175").Formula = "=""09-09-2021"""
$ws.Range("A174:A175").Copy()
$ws.Range("A174:A175").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B174").Value = 1.5
$ws.Range("C174").Value = 1.75
$ws.Range("D174").Value = 1.25

$ws.Range("B175").Value = 1.5
$ws.Range("C175").Value = 1.75
$ws.Range("D175").Value = 1.25
